# ---------------------------------------------------------------------------
# 18/09/2017 CHITRA MAMATHA CHICK IN
#
# 1) The existing "Sun Sep 16 13:21:13 PDT 2017" timestamp line was split
#    across two runs ("Sun Sep 16" + " 13:21:13 PDT 2017"); collapse it to
#    a single run with the full text.
# 2) Append a brand-new purchase-details entry ("Mon Sep 17 14:03:53 PDT
#    2017", Person Name - BM H, Item Name - CARROT, ...) right after the
#    "Amount balance - 339892.0" paragraph that currently ends the document
#    body content.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Change 1: merge the two timestamp runs into one -----------------------
$d.Content.Find.Execute(
    "Sun Sep 16 13:21:13 PDT 2017", $false, $false, $false, $false, $false,
    $true, 1, $false, "Sun Sep 16 13:21:13 PDT 2017", 2) | Out-Null

# --- Change 2: insert the new "Mon Sep 17" entry ----------------------------
# Locate the paragraph that holds "Amount balance ... - 339892.0" (the last
# line of the "Sun Sep 16" entry); the new block goes right after it.
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs($i).Range.Text
    if ($paraText -like "Amount balance*339892.0*") {
        $anchorPara = $d.Paragraphs($i)
        break
    }
}
if ($anchorPara -eq $null) {
    throw "Could not locate the 'Amount balance - 339892.0' anchor paragraph"
}

$insertPoint = $d.Range($anchorPara.Range.End, $anchorPara.Range.End)

# WordprocessingML fragment for the new paragraphs: an empty bold paragraph,
# the "Mon Sep 17 14:03:53 PDT 2017" line, the full purchase-details block
# for CARROT (BM H), ending in "Amount balance - 343292.0", then a blank
# paragraph and a blank bold paragraph (mirrors the previous entry's
# trailing spacer lines).
$rPrPlain = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr>'
$rPrBold  = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/></w:rPr>'
$pPrPlain = '<w:pPr><w:pStyle w:val="PlainText"/>' + $rPrPlain + '</w:pPr>'
$pPrBold  = '<w:pPr><w:pStyle w:val="PlainText"/>' + $rPrBold + '</w:pPr>'

$parts = @(
    "<w:p>$pPrBold</w:p>"

    "<w:p>$pPrPlain" +
        "<w:r>$rPrPlain<w:t>Mon Sep 17</w:t></w:r>" +
        "<w:r>$rPrPlain<w:t xml:space=`"preserve`"> 14:03:53 PDT 2017</w:t></w:r>" +
    "</w:p>"

    "<w:p>$pPrPlain" +
        "<w:r>$rPrPlain<w:t>Person Name</w:t></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/><w:t>- BM H</w:t></w:r>" +
    "</w:p>"

    "<w:p>$pPrPlain" +
        "<w:r>$rPrPlain<w:t>---------------------------------------------------------------</w:t></w:r>" +
    "</w:p>"

    "<w:p>$pPrPlain" +
        "<w:r>$rPrPlain<w:t>Item Name</w:t></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/><w:t>- CARROT</w:t></w:r>" +
    "</w:p>"

    "<w:p>$pPrPlain" +
        "<w:r>$rPrPlain<w:t>Number of Pockets</w:t></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/><w:t>- 2</w:t></w:r>" +
    "</w:p>"

    "<w:p>$pPrPlain" +
        "<w:r>$rPrPlain<w:t>Number of KGs</w:t></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/><w:t>- 169</w:t></w:r>" +
    "</w:p>"

    "<w:p>$pPrPlain" +
        "<w:r>$rPrPlain<w:t>Rate</w:t></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/><w:t>- 20</w:t></w:r>" +
    "</w:p>"

    "<w:p>$pPrPlain" +
        "<w:r>$rPrPlain<w:t>Transport &amp; Miscellaneous</w:t></w:r>" +
        "<w:r>$rPrPlain<w:tab/><w:t>- 20</w:t></w:r>" +
    "</w:p>"

    "<w:p>$pPrPlain" +
        "<w:r>$rPrPlain<w:t>Total Price</w:t></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/></w:r>" +
        "<w:r>$rPrPlain<w:tab/><w:t>- 3400.0</w:t></w:r>" +
    "</w:p>"

    "<w:p>$pPrBold" +
        "<w:r>$rPrBold<w:t>Amount balance</w:t></w:r>" +
        "<w:r>$rPrBold<w:tab/></w:r>" +
        "<w:r>$rPrBold<w:tab/></w:r>" +
        "<w:r>$rPrBold<w:tab/><w:t>- 343292.0</w:t></w:r>" +
    "</w:p>"

    "<w:p>$pPrPlain</w:p>"

    "<w:p>$pPrBold</w:p>"
)

$fragment = [string]::Join("", $parts)
$insertPoint.InsertXML($fragment) | Out-Null
